{"js": "// Replace the date heading and every arithmetic expression in the\n// practice-sheet table with the next day's worksheet values.\n//\n// The mapping below is the ordered list of [oldText, newText] pairs,\n// matching the document order of: title paragraph, then every table\n// cell's paragraph in row-major order (20 rows x 5 columns).\nconst replacements = [\n  [\"2024-08-25 Sunday\", \"2024-08-26 Monday\"],\n  [\"5+38=\", \"16+50=\"],\n  [\"77-67=\", \"97-36=\"],\n  [\"90+6=\", \"78-8=\"],\n  [\"19+38=\", \"65+2=\"],\n  [\"89-87=\", \"81-45=\"],\n  [\"27+9=\", \"94-70=\"],\n  [\"55-34=\", \"14+8=\"],\n  [\"7+18=\", \"32-12=\"],\n  [\"77-22=\", \"97-24=\"],\n  [\"41+29=\", \"47-15=\"],\n  [\"70+14=\", \"72-6=\"],\n  [\"64-61=\", \"85+12=\"],\n  [\"87-33=\", \"35+21=\"],\n  [\"39-14=\", \"66-51=\"],\n  [\"7+0=\", \"96-41=\"],\n  [\"13+85=\", \"9+50=\"],\n  [\"34+8=\", \"83-56=\"],\n  [\"67+2=\", \"73-66=\"],\n  [\"3-0=\", \"61+36=\"],\n  [\"73-70=\", \"68-59=\"],\n  [\"66+32=\", \"61+20=\"],\n  [\"6+68=\", \"98-89=\"],\n  [\"1+92=\", \"19-11=\"],\n  [\"42+39=\", \"24+66=\"],\n  [\"11+80=\", \"90-17=\"],\n  [\"99-53=\", \"62-17=\"],\n  [\"81-64=\", \"65+13=\"],\n  [\"68+15=\", \"92-18=\"],\n  [\"39-23=\", \"6+11=\"],\n  [\"57-29=\", \"88-82=\"],\n  [\"66-30=\", \"85-15=\"],\n  [\"68+12=\", \"56+31=\"],\n  [\"61-7=\", \"46+20=\"],\n  [\"70-66=\", \"14+41=\"],\n  [\"40+13=\", \"44+0=\"],\n  [\"98-59=\", \"43+30=\"],\n  [\"7+38=\", \"24+69=\"],\n  [\"60-9=\", \"90-76=\"],\n  [\"14+60=\", \"65+33=\"],\n  [\"78-72=\", \"22+56=\"],\n  [\"12+9=\", \"1+81=\"],\n  [\"63-32=\", \"97-58=\"],\n  [\"50-3=\", \"5+83=\"],\n  [\"55-10=\", \"60-57=\"],\n  [\"42+15=\", \"51-37=\"],\n  [\"37-34=\", \"96-71=\"],\n  [\"58+20=\", \"23-2=\"],\n  [\"89-82=\", \"65-21=\"],\n  [\"84-3=\", \"84-79=\"],\n  [\"4+92=\", \"63-62=\"],\n  [\"96-35=\", \"86+3=\"],\n  [\"89+6=\", \"47+32=\"],\n  [\"56+5=\", \"3+68=\"],\n  [\"72-45=\", \"89-11=\"],\n  [\"72+10=\", \"97-28=\"],\n  [\"67-10=\", \"81-29=\"],\n  [\"68-56=\", \"39+45=\"],\n  [\"68-45=\", \"71-60=\"],\n  [\"14+79=\", \"68+22=\"],\n  [\"48+24=\", \"12+12=\"],\n  [\"18+7=\", \"29+69=\"],\n  [\"64-19=\", \"22-7=\"],\n  [\"76+3=\", \"32-29=\"],\n  [\"69-6=\", \"97-58=\"],\n  [\"27-26=\", \"18+8=\"],\n  [\"74-2=\", \"22+8=\"],\n  [\"86-49=\", \"81-26=\"],\n  [\"90-54=\", \"43-12=\"],\n  [\"3+40=\", \"60-14=\"],\n  [\"18+19=\", \"0+41=\"],\n  [\"44+32=\", \"38+29=\"],\n  [\"59+36=\", \"55-19=\"],\n  [\"38-12=\", \"92-53=\"],\n  [\"96-86=\", \"71-25=\"],\n  [\"30-23=\", \"89-18=\"],\n  [\"33-20=\", \"63-63=\"],\n  [\"2+15=\", \"30+10=\"],\n  [\"26+68=\", \"27-23=\"],\n  [\"37-30=\", \"99-83=\"],\n  [\"0+33=\", \"16+42=\"],\n  [\"22+6=\", \"97-7=\"],\n  [\"11+44=\", \"67-32=\"],\n  [\"55+36=\", \"15+5=\"],\n  [\"56+7=\", \"9+0=\"],\n  [\"44-14=\", \"13+38=\"],\n  [\"23+15=\", \"62+3=\"],\n  [\"76-54=\", \"54+1=\"],\n  [\"84+7=\", \"60+31=\"],\n  [\"63+23=\", \"1+82=\"],\n  [\"19+52=\", \"70+11=\"],\n  [\"88-74=\", \"16+22=\"],\n  [\"66+11=\", \"58+22=\"],\n  [\"56-8=\", \"0+82=\"],\n  [\"29+43=\", \"28+52=\"],\n  [\"32+38=\", \"19+45=\"],\n  [\"60-5=\", \"78+13=\"],\n  [\"58+4=\", \"27+68=\"],\n  [\"92-22=\", \"4+3=\"],\n  [\"95-77=\", \"5+45=\"],\n  [\"29-15=\", \"76-53=\"]\n];\n\n// `body.paragraphs` walks the whole document in order and (for this\n// document) yields exactly: the title paragraph followed by every\n// table-cell paragraph in row-major order, one text run each -- lining\n// up 1:1 with `replacements` above.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nif (items.length !== replacements.length) {\n  throw new Error(\n    `Expected ${replacements.length} paragraphs, found ${items.length}`\n  );\n}\n\nfor (let i = 0; i < items.length; i++) {\n  const [oldText, newText] = replacements[i];\n  const para = items[i];\n  if (para.text !== oldText) {\n    throw new Error(\n      `Paragraph ${i}: expected \"${oldText}\" but found \"${para.text}\"`\n    );\n  }\n  // insertText(..., replace) swaps the paragraph's text while keeping\n  // the existing run formatting (rFonts/sz/etc.) intact.\n  para.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date heading and every arithmetic expression in\n# the practice table to the next day's values.\n$d = $word.ActiveDocument\n\n# --- 1. Title paragraph (date line) ---\n$titleOld = \"2024-08-25 Sunday\"\n$titleNew = \"2024-08-26 Monday\"\n$p1 = $d.Paragraphs.First\nif ($p1.Range.Text.TrimEnd([char]13, [char]7) -ne $titleOld) {\n    throw \"Title paragraph text mismatch: found '$($p1.Range.Text)'\"\n}\n$p1.Range.Text = $titleNew\n\n# --- 2. Table cells, row-major order (20 rows x 5 columns) ---\n$cellPairs = @(\n    @(\"5+38=\", \"16+50=\"),\n    @(\"77-67=\", \"97-36=\"),\n    @(\"90+6=\", \"78-8=\"),\n    @(\"19+38=\", \"65+2=\"),\n    @(\"89-87=\", \"81-45=\"),\n    @(\"27+9=\", \"94-70=\"),\n    @(\"55-34=\", \"14+8=\"),\n    @(\"7+18=\", \"32-12=\"),\n    @(\"77-22=\", \"97-24=\"),\n    @(\"41+29=\", \"47-15=\"),\n    @(\"70+14=\", \"72-6=\"),\n    @(\"64-61=\", \"85+12=\"),\n    @(\"87-33=\", \"35+21=\"),\n    @(\"39-14=\", \"66-51=\"),\n    @(\"7+0=\", \"96-41=\"),\n    @(\"13+85=\", \"9+50=\"),\n    @(\"34+8=\", \"83-56=\"),\n    @(\"67+2=\", \"73-66=\"),\n    @(\"3-0=\", \"61+36=\"),\n    @(\"73-70=\", \"68-59=\"),\n    @(\"66+32=\", \"61+20=\"),\n    @(\"6+68=\", \"98-89=\"),\n    @(\"1+92=\", \"19-11=\"),\n    @(\"42+39=\", \"24+66=\"),\n    @(\"11+80=\", \"90-17=\"),\n    @(\"99-53=\", \"62-17=\"),\n    @(\"81-64=\", \"65+13=\"),\n    @(\"68+15=\", \"92-18=\"),\n    @(\"39-23=\", \"6+11=\"),\n    @(\"57-29=\", \"88-82=\"),\n    @(\"66-30=\", \"85-15=\"),\n    @(\"68+12=\", \"56+31=\"),\n    @(\"61-7=\", \"46+20=\"),\n    @(\"70-66=\", \"14+41=\"),\n    @(\"40+13=\", \"44+0=\"),\n    @(\"98-59=\", \"43+30=\"),\n    @(\"7+38=\", \"24+69=\"),\n    @(\"60-9=\", \"90-76=\"),\n    @(\"14+60=\", \"65+33=\"),\n    @(\"78-72=\", \"22+56=\"),\n    @(\"12+9=\", \"1+81=\"),\n    @(\"63-32=\", \"97-58=\"),\n    @(\"50-3=\", \"5+83=\"),\n    @(\"55-10=\", \"60-57=\"),\n    @(\"42+15=\", \"51-37=\"),\n    @(\"37-34=\", \"96-71=\"),\n    @(\"58+20=\", \"23-2=\"),\n    @(\"89-82=\", \"65-21=\"),\n    @(\"84-3=\", \"84-79=\"),\n    @(\"4+92=\", \"63-62=\"),\n    @(\"96-35=\", \"86+3=\"),\n    @(\"89+6=\", \"47+32=\"),\n    @(\"56+5=\", \"3+68=\"),\n    @(\"72-45=\", \"89-11=\"),\n    @(\"72+10=\", \"97-28=\"),\n    @(\"67-10=\", \"81-29=\"),\n    @(\"68-56=\", \"39+45=\"),\n    @(\"68-45=\", \"71-60=\"),\n    @(\"14+79=\", \"68+22=\"),\n    @(\"48+24=\", \"12+12=\"),\n    @(\"18+7=\", \"29+69=\"),\n    @(\"64-19=\", \"22-7=\"),\n    @(\"76+3=\", \"32-29=\"),\n    @(\"69-6=\", \"97-58=\"),\n    @(\"27-26=\", \"18+8=\"),\n    @(\"74-2=\", \"22+8=\"),\n    @(\"86-49=\", \"81-26=\"),\n    @(\"90-54=\", \"43-12=\"),\n    @(\"3+40=\", \"60-14=\"),\n    @(\"18+19=\", \"0+41=\"),\n    @(\"44+32=\", \"38+29=\"),\n    @(\"59+36=\", \"55-19=\"),\n    @(\"38-12=\", \"92-53=\"),\n    @(\"96-86=\", \"71-25=\"),\n    @(\"30-23=\", \"89-18=\"),\n    @(\"33-20=\", \"63-63=\"),\n    @(\"2+15=\", \"30+10=\"),\n    @(\"26+68=\", \"27-23=\"),\n    @(\"37-30=\", \"99-83=\"),\n    @(\"0+33=\", \"16+42=\"),\n    @(\"22+6=\", \"97-7=\"),\n    @(\"11+44=\", \"67-32=\"),\n    @(\"55+36=\", \"15+5=\"),\n    @(\"56+7=\", \"9+0=\"),\n    @(\"44-14=\", \"13+38=\"),\n    @(\"23+15=\", \"62+3=\"),\n    @(\"76-54=\", \"54+1=\"),\n    @(\"84+7=\", \"60+31=\"),\n    @(\"63+23=\", \"1+82=\"),\n    @(\"19+52=\", \"70+11=\"),\n    @(\"88-74=\", \"16+22=\"),\n    @(\"66+11=\", \"58+22=\"),\n    @(\"56-8=\", \"0+82=\"),\n    @(\"29+43=\", \"28+52=\"),\n    @(\"32+38=\", \"19+45=\"),\n    @(\"60-5=\", \"78+13=\"),\n    @(\"58+4=\", \"27+68=\"),\n    @(\"92-22=\", \"4+3=\"),\n    @(\"95-77=\", \"5+45=\"),\n    @(\"29-15=\", \"76-53=\")\n)\n\n$t = $d.Tables.Item(1)\n$cols = $t.Columns.Count\n$rows = $t.Rows.Count\n\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $pair = $cellPairs[$i]\n        $oldText = $pair[0]\n        $newText = $pair[1]\n\n        $cell = $t.Cell($r, $c)\n        $actual = $cell.Range.Text.TrimEnd([char]13, [char]7)\n        if ($actual -ne $oldText) {\n            throw \"Cell ($r,$c) text mismatch: expected '$oldText' but found '$actual'\"\n        }\n        $cell.Range.Text = $newText\n\n        $i++\n    }\n}\n"}
